$d = $word.ActiveDocument

# 1) Remove the "both" (justify) paragraph alignment from the paragraph that
#    contains "We seek to create a catalog..." (the pPr's <w:jc w:val="both"/>).
$alignRng = $d.Content
$alignRng.Find.Execute("We seek to create a catalog", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($alignRng.Find.Found) {
    $alignRng.Paragraphs.First.Format.Alignment = 0
}

# 2) Move the "_GoBack" bookmark from right before "positives, utilizing"
#    to right after the word "is" in "...campaigns 4 and 5, which is expected...".
#    Adding a bookmark named "_GoBack" relocates any existing one of that name.
$isRng = $d.Content
$isRng.Find.Execute("is", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($isRng.Find.Found) {
    $isRng.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $isRng)
}

# 3) Now that the bookmark no longer sits between "false " and "positives,
#    utilizing..." re-run a Find/Replace over that text so Word coalesces the
#    (now contiguous, identically-formatted) runs into a single run.
$mergeRng = $d.Content
$mergeRng.Find.Execute("false positives, utilizing", $false, $false, $false, $false, $false, $true, 1, $false, "false positives, utilizing", 2)
